$d = $word.ActiveDocument
$t = $d.Tables.Item(1)

# ---------------------------------------------------------------------------
# Helper pattern: split a date-like "DD.MM" run into three runs "D","D",".MM"
# (same total text, or a changed digit) while keeping the identical run
# formatting. Setting Bold on/off on each sub-range after the text has been
# written forces the engine to keep them as separate <w:r> elements instead
# of silently re-merging adjacent runs that share identical rPr.
# ---------------------------------------------------------------------------

function Split-DateCell($row, $col, $d0, $d1, $tail) {
    $cell = $t.Cell($row, $col)
    $start = $cell.Range.Start

    $rA = $d.Range($start, $start + 1)
    $rA.Text = $d0

    $rB = $d.Range($start + 1, $start + 2)
    $rB.Text = $d1

    $rC = $d.Range($start + 2, $start + 2 + $tail.Length)
    $rC.Text = $tail

    # Re-select the (possibly re-seated) sub ranges and pin the run
    # boundaries by toggling a no-op character format.
    $rA2 = $d.Range($start, $start + 1)
    $rA2.Bold = 1
    $rA2.Bold = 0

    $rB2 = $d.Range($start + 1, $start + 2)
    $rB2.Bold = 1
    $rB2.Bold = 0

    $rC2 = $d.Range($start + 2, $start + 2 + $tail.Length)
    $rC2.Bold = 1
    $rC2.Bold = 0
}

# Row 43 (ЛР16 / "Заняття" column): 03.05 -> 04.05
Split-DateCell 43 2 "0" "4" ".05"

# Row 45 (ЛР17 / "Заняття" column): 05.05 -> 05.05 (re-split into 3 runs)
Split-DateCell 45 2 "0" "5" ".05"

# Row 47 (Л25 / "Заняття" column): 06.05 -> 07.05
Split-DateCell 47 2 "0" "7" ".05"

# ---------------------------------------------------------------------------
# Row 50 (ЛР19 "Тема" column): drop the trailing, italic/red
# " Підсумкове тематичне тестування" phrase down to a single space while
# keeping that run's formatting untouched.
# ---------------------------------------------------------------------------
$cell50 = $t.Cell(50, 5)
$found = $cell50.Range.Find.Execute(" Підсумкове тематичне тестування", $true, $false, $false, $false, $false, $true, 1, $false, " ", 2)

# ---------------------------------------------------------------------------
# Row 52 (Л27 "Тема" column): append a new italic/red
# " Підсумкове тематичне тестування" run right after the trailing "." run.
# Inserting off of that run's own Range makes the new text inherit its
# formatting (Times New Roman incl. w:cs, sz/szCs 28, italic); we then only
# need to flip the new span's color to red.
# ---------------------------------------------------------------------------
$cell52 = $t.Cell(52, 5)
$full52 = $cell52.Range.Text
$start52 = $cell52.Range.Start
$textLen52 = $full52.Length - 2

$dotRange = $d.Range($start52 + $textLen52 - 1, $start52 + $textLen52)
$dotRange.InsertAfter(" Підсумкове тематичне тестування")

$cell52b = $t.Cell(52, 5)
$full52b = $cell52b.Range.Text
$start52b = $cell52b.Range.Start
$textLen52b = $full52b.Length - 2
$suffixLen = " Підсумкове тематичне тестування".Length

$newRange = $d.Range($start52b + $textLen52b - $suffixLen, $start52b + $textLen52b)
$newRange.Font.Color = 255
